$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.002837
$ws.Range("H2").Value = 0.008510999999999999
$ws.Range("I2").Value = 0.00007108247730492929
$ws.Range("J2").Value = 0.00007108247730492929
$ws.Range("M2").Value = 2.685464
$ws.Range("N2").Value = 8.056392000000001
$ws.Range("O2").Value = 0.06676031826184478
$ws.Range("P2").Value = 0.06676031826184478
$ws.Range("Q2").Value = 0.007618661367999999
$ws.Range("R2").Value = 0.068567952312
$ws.Range("S2").Value = 0.000004745488807717438
$ws.Range("T2").Value = 0.000004745488807717438
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.002837
$ws.Range("H3").Value = 0.008510999999999999
$ws.Range("I3").Value = 0.00007108247730492929
$ws.Range("J3").Value = 0.00007108247730492929
$ws.Range("O3").Value = 0.02342101692711854
$ws.Range("P3").Value = 0.02342101692711854
$ws.Range("Q3").Value = 0.002672797277
$ws.Range("R3").Value = 0.024055175493
$ws.Range("S3").Value = 0.000001664823904180268
$ws.Range("T3").Value = 0.000001664823904180268
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.002837
$ws.Range("H4").Value = 0.008510999999999999
$ws.Range("I4").Value = 0.00007108247730492929
$ws.Range("J4").Value = 0.00007108247730492929
$ws.Range("M4").Value = 35.399925
$ws.Range("N4").Value = 106.199775
$ws.Range("O4").Value = 0.8800379597140142
$ws.Range("P4").Value = 0.8800379597140142
$ws.Range("Q4").Value = 0.100429587225
$ws.Range("R4").Value = 0.903866285025
$ws.Range("S4").Value = 0.00006255527829884768
$ws.Range("T4").Value = 0.00006255527829884768
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.002837
$ws.Range("H5").Value = 0.008510999999999999
$ws.Range("I5").Value = 0.00007108247730492929
$ws.Range("J5").Value = 0.00007108247730492929
$ws.Range("M5").Value = 1.197942333333333
$ws.Range("N5").Value = 3.593827
$ws.Range("O5").Value = 0.02978070509702244
$ws.Range("P5").Value = 0.02978070509702244
$ws.Range("Q5").Value = 0.003398562399666666
$ws.Range("R5").Value = 0.030587061597
$ws.Range("S5").Value = 0.00000211688629418389
$ws.Range("T5").Value = 0.00000211688629418389
$ws.Range("I6").Value = 0.3776915775490952
$ws.Range("J6").Value = 0.3776915775490952
$ws.Range("M6").Value = 2.685464
$ws.Range("N6").Value = 8.056392000000001
$ws.Range("O6").Value = 0.06676031826184478
$ws.Range("P6").Value = 0.06676031826184478
$ws.Range("Q6").Value = 40.48120352570667
$ws.Range("R6").Value = 364.33083173136
$ws.Range("S6").Value = 0.02521480992199582
$ws.Range("T6").Value = 0.02521480992199582
$ws.Range("I7").Value = 0.3776915775490952
$ws.Range("J7").Value = 0.3776915775490952
$ws.Range("O7").Value = 0.02342101692711854
$ws.Range("P7").Value = 0.02342101692711854
$ws.Range("S7").Value = 0.008845920831007464
$ws.Range("T7").Value = 0.008845920831007464
$ws.Range("I8").Value = 0.3776915775490952
$ws.Range("J8").Value = 0.3776915775490952
$ws.Range("M8").Value = 35.399925
$ws.Range("N8").Value = 106.199775
$ws.Range("O8").Value = 0.8800379597140142
$ws.Range("P8").Value = 0.8800379597140142
$ws.Range("Q8").Value = 533.6253134355001
$ws.Range("R8").Value = 4802.6278209195
$ws.Range("S8").Value = 0.3323829253074732
$ws.Range("T8").Value = 0.3323829253074732
$ws.Range("I9").Value = 0.3776915775490952
$ws.Range("J9").Value = 0.3776915775490952
$ws.Range("M9").Value = 1.197942333333333
$ws.Range("N9").Value = 3.593827
$ws.Range("O9").Value = 0.02978070509702244
$ws.Range("P9").Value = 0.02978070509702244
$ws.Range("Q9").Value = 18.05801433485111
$ws.Range("R9").Value = 162.52212901366
$ws.Range("S9").Value = 0.01124792148861879
$ws.Range("T9").Value = 0.01124792148861879
$ws.Range("G10").Value = 1.581618666666667
$ws.Range("H10").Value = 4.744856
$ws.Range("I10").Value = 0.03962825977384063
$ws.Range("J10").Value = 0.03962825977384063
$ws.Range("M10").Value = 2.685464
$ws.Range("N10").Value = 8.056392000000001
$ws.Range("O10").Value = 0.06676031826184478
$ws.Range("P10").Value = 0.06676031826184478
$ws.Range("Q10").Value = 4.247379991061334
$ws.Range("R10").Value = 38.22641991955201
$ws.Range("S10").Value = 0.002645595234664662
$ws.Range("T10").Value = 0.002645595234664662
$ws.Range("G11").Value = 1.581618666666667
$ws.Range("H11").Value = 4.744856
$ws.Range("I11").Value = 0.03962825977384063
$ws.Range("J11").Value = 0.03962825977384063
$ws.Range("O11").Value = 0.02342101692711854
$ws.Range("P11").Value = 0.02342101692711854
$ws.Range("Q11").Value = 1.490076159858667
$ws.Range("R11").Value = 13.410685438728
$ws.Range("S11").Value = 0.0009281341429553721
$ws.Range("T11").Value = 0.0009281341429553721
$ws.Range("G12").Value = 1.581618666666667
$ws.Range("H12").Value = 4.744856
$ws.Range("I12").Value = 0.03962825977384063
$ws.Range("J12").Value = 0.03962825977384063
$ws.Range("M12").Value = 35.399925
$ws.Range("N12").Value = 106.199775
$ws.Range("O12").Value = 0.8800379597140142
$ws.Range("P12").Value = 0.8800379597140142
$ws.Range("Q12").Value = 55.9891821786
$ws.Range("R12").Value = 503.9026396074
$ws.Range("S12").Value = 0.03487437287838765
$ws.Range("T12").Value = 0.03487437287838765
$ws.Range("G13").Value = 1.581618666666667
$ws.Range("H13").Value = 4.744856
$ws.Range("I13").Value = 0.03962825977384063
$ws.Range("J13").Value = 0.03962825977384063
$ws.Range("M13").Value = 1.197942333333333
$ws.Range("N13").Value = 3.593827
$ws.Range("O13").Value = 0.02978070509702244
$ws.Range("P13").Value = 0.02978070509702244
$ws.Range("Q13").Value = 1.894687955990223
$ws.Range("R13").Value = 17.052191603912
$ws.Range("S13").Value = 0.001180157517832945
$ws.Range("T13").Value = 0.001180157517832945
$ws.Range("G14").Value = 23.25273433333334
$ws.Range("H14").Value = 69.75820300000001
$ws.Range("I14").Value = 0.5826090801997593
$ws.Range("J14").Value = 0.5826090801997593
$ws.Range("M14").Value = 2.685464
$ws.Range("N14").Value = 8.056392000000001
$ws.Range("O14").Value = 0.06676031826184478
$ws.Range("P14").Value = 0.06676031826184478
$ws.Range("Q14").Value = 62.44438095373068
$ws.Range("R14").Value = 561.9994285835761
$ws.Range("S14").Value = 0.03889516761637658
$ws.Range("T14").Value = 0.03889516761637658
$ws.Range("G15").Value = 23.25273433333334
$ws.Range("H15").Value = 69.75820300000001
$ws.Range("I15").Value = 0.5826090801997593
$ws.Range("J15").Value = 0.5826090801997593
$ws.Range("O15").Value = 0.02342101692711854
$ws.Range("P15").Value = 0.02342101692711854
$ws.Range("Q15").Value = 21.90688932285434
$ws.Range("R15").Value = 197.162003905689
$ws.Range("S15").Value = 0.01364529712925153
$ws.Range("T15").Value = 0.01364529712925153
$ws.Range("G16").Value = 23.25273433333334
$ws.Range("H16").Value = 69.75820300000001
$ws.Range("I16").Value = 0.5826090801997593
$ws.Range("J16").Value = 0.5826090801997593
$ws.Range("M16").Value = 35.399925
$ws.Range("N16").Value = 106.199775
$ws.Range("O16").Value = 0.8800379597140142
$ws.Range("P16").Value = 0.8800379597140142
$ws.Range("Q16").Value = 823.1450514449252
$ws.Range("R16").Value = 7408.305463004326
$ws.Range("S16").Value = 0.5127181062498547
$ws.Range("T16").Value = 0.5127181062498547
$ws.Range("G17").Value = 23.25273433333334
$ws.Range("H17").Value = 69.75820300000001
$ws.Range("I17").Value = 0.5826090801997593
$ws.Range("J17").Value = 0.5826090801997593
$ws.Range("M17").Value = 1.197942333333333
$ws.Range("N17").Value = 3.593827
$ws.Range("O17").Value = 0.02978070509702244
$ws.Range("P17").Value = 0.02978070509702244
$ws.Range("Q17").Value = 27.85543482365345
$ws.Range("R17").Value = 250.698913412881
$ws.Range("S17").Value = 0.01735050920427653
$ws.Range("T17").Value = 0.01735050920427653
